$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$s.Shapes.Item("Rectangle 1").Delete()
